# fix: alterar python version para 3.11.5
# Update absenteeism data rows 2-11 with new values as per upstream diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 61958
$ws.Range("B2").Value = "Sr. João Vitor Barros"
$ws.Range("C2").Value = "Atendimento ao Cliente"
$ws.Range("D2").Value = "Outros"
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 45103
$ws.Range("G2").Value = 11873.74

# Row 3
$ws.Range("A3").Value = 44807
$ws.Range("B3").Value = "Samuel Ribeiro"
$ws.Range("D3").Value = "Doença"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 45091
$ws.Range("G3").Value = 10271.22

# Row 4
$ws.Range("A4").Value = 62826
$ws.Range("B4").Value = "Pedro Lucas Azevedo"
$ws.Range("C4").Value = "Engenharia"
$ws.Range("D4").Value = "Problemas pessoais"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 45097
$ws.Range("G4").Value = 4164.59

# Row 5
$ws.Range("A5").Value = 55791
$ws.Range("B5").Value = "Davi Cunha"
$ws.Range("C5").Value = "Marketing"
$ws.Range("F5").Value = 45089
$ws.Range("G5").Value = 4485.45

# Row 6
$ws.Range("A6").Value = 55703
$ws.Range("B6").Value = "Bruna Porto"
$ws.Range("C6").Value = "Marketing"
$ws.Range("F6").Value = 45103
$ws.Range("G6").Value = 5870.6

# Row 7
$ws.Range("A7").Value = 8007
$ws.Range("B7").Value = "Sofia Barbosa"
$ws.Range("C7").Value = "Atendimento ao Cliente"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 45094
$ws.Range("G7").Value = 10371.15

# Row 8
$ws.Range("A8").Value = 25750
$ws.Range("B8").Value = "Agatha Cardoso"
$ws.Range("C8").Value = "Marketing"
$ws.Range("D8").Value = "Consulta médica"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 45102
$ws.Range("G8").Value = 9323.99

# Row 9
$ws.Range("A9").Value = 69783
$ws.Range("B9").Value = "Ana Laura Aragão"
$ws.Range("C9").Value = "Jurídico"
$ws.Range("D9").Value = "Viagem de negócios"
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = 45099
$ws.Range("G9").Value = 9779

# Row 10
$ws.Range("A10").Value = 50311
$ws.Range("B10").Value = "Rafaela Novaes"
$ws.Range("C10").Value = "Recursos Humanos"
$ws.Range("D10").Value = "Viagem de negócios"
$ws.Range("E10").Value = 8
$ws.Range("F10").Value = 45096
$ws.Range("G10").Value = 9288.53

# Row 11
$ws.Range("A11").Value = 26622
$ws.Range("B11").Value = "Sra. Sofia da Mota"
$ws.Range("C11").Value = "TI"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 45106
$ws.Range("G11").Value = 9742.13

$wb.Save()
